# Updated symbol list with GitHub Actions refresh of crypto price/volume data.
# Values are entered with a leading apostrophe so Excel stores them as literal
# text (matching the workbook's existing inlineStr cells) instead of re-typing
# them as numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.60"
$ws.Range("E2").Value = "'1.21%"

$ws.Range("D3").Value = "'29.32"
$ws.Range("E3").Value = "'-1.50%"

$ws.Range("D4").Value = "'5.158"
$ws.Range("E4").Value = "'0.69%"

$ws.Range("D5").Value = "'0.05764"
$ws.Range("E5").Value = "'1.99%"

$ws.Range("D6").Value = "'6.600"
$ws.Range("E6").Value = "'1.61%"

$ws.Range("D7").Value = "'0.8592"
$ws.Range("E7").Value = "'3.93%"

$ws.Range("D8").Value = "'0.8606"
$ws.Range("E8").Value = "'-0.22%"

$ws.Range("D9").Value = "'0.1365"
$ws.Range("E9").Value = "'2.56%"

$ws.Range("D10").Value = "'0.07023"
$ws.Range("E10").Value = "'1.55%"

$ws.Range("D11").Value = "'0.03026"
$ws.Range("E11").Value = "'5.91%"

$ws.Range("D12").Value = "'0.09364"
$ws.Range("E12").Value = "'-0.22%"

$ws.Range("D13").Value = "'0.001523"
$ws.Range("E13").Value = "'0.38%"

$ws.Range("D14").Value = "'0.0006029"
$ws.Range("E14").Value = "'-94.01%"

$ws.Range("D15").Value = "'0.006061"
$ws.Range("E15").Value = "'-0.63%"

$ws.Range("E16").Value = "'-0.80%"

$ws.Range("D17").Value = "'3.146"
$ws.Range("E17").Value = "'4.41%"

$ws.Range("D18").Value = "'2.154"
$ws.Range("E18").Value = "'-2.74%"

$ws.Range("D19").Value = "'0.3201"
$ws.Range("E19").Value = "'1.67%"

$ws.Range("D20").Value = "'0.03311"
$ws.Range("E20").Value = "'1.92%"

$ws.Range("D21").Value = "'0.1282"
$ws.Range("E21").Value = "'-0.96%"

$ws.Range("D22").Value = "'3.318"
$ws.Range("E22").Value = "'-8.25%"

$ws.Range("D23").Value = "'0.04145"
$ws.Range("E23").Value = "'-0.29%"

$ws.Range("D24").Value = "'0.1401"
$ws.Range("E24").Value = "'1.97%"

$ws.Range("D25").Value = "'0.001226"
$ws.Range("E25").Value = "'1.37%"

$ws.Range("D26").Value = "'0.004134"
$ws.Range("E26").Value = "'-6.97%"

$ws.Range("D27").Value = "'0.0001210"
$ws.Range("E27").Value = "'2.53%"

$ws.Range("E28").Value = "'3.22%"

$ws.Range("D40").Value = "'0.03732"
$ws.Range("E40").Value = "'0.68%"

$ws.Range("D41").Value = "'0.005878"
$ws.Range("E41").Value = "'2.25%"

$ws.Range("E42").Value = "'1.53%"

$ws.Range("D43").Value = "'0.002199"
$ws.Range("E43").Value = "'-4.84%"

$ws.Range("D44").Value = "'0.008392"
$ws.Range("E44").Value = "'-13.55%"

$ws.Range("D45").Value = "'0.00005292"
$ws.Range("E45").Value = "'3.66%"

$ws.Range("E46").Value = "'-0.02%"

$ws.Range("D47").Value = "'0.05799"
$ws.Range("E47").Value = "'-44.75%"

$ws.Range("D48").Value = "'0.002445"
$ws.Range("E48").Value = "'-15.99%"

$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'-0.02%"

$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'-0.02%"
